$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.064.61'
$ws.Range('E2').Value = '  -7.80%  '
$ws.Range('D3').Value = '3.259.47'
$ws.Range('E3').Value = '  -9.66%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '176.28'
$ws.Range('E5').Value = '  -14.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '509.82'
$ws.Range('E6').Value = '  -10.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.586'
$ws.Range('E7').Value = '  -4.49%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '3.253.97'
$ws.Range('E9').Value = '  -9.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.612'
$ws.Range('E10').Value = '  -10.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.85'
$ws.Range('E11').Value = '  -11.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.129'
$ws.Range('E12').Value = '  -13.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('E13').Value = '  -11.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.04'
$ws.Range('E14').Value = '  -12.36%  '
$ws.Range('D15').Value = '3.747.12'
$ws.Range('E15').Value = '  -10.45%  '
$ws.Range('E16').Value = '  -6.37%  '
$ws.Range('D17').Value = '3.237.84'
$ws.Range('E17').Value = '  -10.73%  '
$ws.Range('D18').Value = '62.866.05'
$ws.Range('E18').Value = '  -7.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.16'
$ws.Range('E19').Value = '  -11.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.82'
$ws.Range('E20').Value = '  -12.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.935'
$ws.Range('E21').Value = '  -12.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '366.23'
$ws.Range('E22').Value = '  -9.83%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.04'
$ws.Range('E23').Value = '  -11.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.84'
$ws.Range('E24').Value = '  -7.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.61'
$ws.Range('E25').Value = '  -14.13%  '
$ws.Range('E26').Value = '  -2.83%  '
$ws.Range('E27').Value = '  -3.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.62'
$ws.Range('E28').Value = '  -10.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.20'
$ws.Range('E29').Value = '  -10.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.20'
$ws.Range('E30').Value = '  -11.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '646.63'
$ws.Range('E31').Value = '  -9.94%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.04'
$ws.Range('E32').Value = '  -11.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.59'
$ws.Range('E33').Value = '  -15.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.05'
$ws.Range('E34').Value = '  -9.53%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.04'
$ws.Range('E35').Value = '  -8.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.103'
$ws.Range('E36').Value = '  -10.24%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.60'
$ws.Range('E38').Value = '  -15.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.377'
$ws.Range('E39').Value = '  -10.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  -0.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.124'
$ws.Range('E41').Value = '  -6.98%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '28.21'
$ws.Range('E42').Value = '  +27.08%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.859.95'
$ws.Range('E43').Value = '  -10.53%  '
$ws.Range('B44').Value = 'PEPE'
$ws.Range('C44').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D44').Value = '0.0₃0646'
$ws.Range('E44').Value = '  -14.46%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.41'
$ws.Range('E45').Value = '  -10.23%  '
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.62'
$ws.Range('E46').Value = '  -20.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.80'
$ws.Range('E47').Value = '  +3.18%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.54'
$ws.Range('E48').Value = '  -8.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0377'
$ws.Range('E49').Value = '  -9.59%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.92'
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.122'
$ws.Range('E51').Value = '  -7.00%  '
